# Edit script: update the corrosion-damage section ("FSec1_corrosion") geometry,
# reassign Beam B4 and Column C4 away from the corroded section back to "FSec1",
# and fix the end-length-offset values for B2/B3 so computed deformations match
# the applied static loads (per commit message).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Frame Sec Def - Steel Pipe": update FSec1_corrosion geometry (row 5)
#    Outside Diameter 293.5 -> 287, Wall Thickness 13 -> 6.5
# ---------------------------------------------------------------------------
$wsPipe = $wb.Worksheets.Item("Frame Sec Def - Steel Pipe")
$wsPipe.Range("D5").Value = 287
$wsPipe.Range("E5").Value = 6.5

# ---------------------------------------------------------------------------
# 2) "Frame Prop - Summary": recalculated section properties for FSec1_corrosion
#    (row 5), consistent with the new pipe geometry above.
# ---------------------------------------------------------------------------
$wsProp = $wb.Worksheets.Item("Frame Prop - Summary")
$wsProp.Range("E5").Value  = 5727.9        # Area
$wsProp.Range("F5").Value  = 112728825.2   # J
$wsProp.Range("G5").Value  = 56364412.6    # I33
$wsProp.Range("H5").Value  = 56364412.6    # I22
$wsProp.Range("I5").Value  = 2865          # As2
$wsProp.Range("J5").Value  = 2865          # As3
$wsProp.Range("K5").Value  = 392783.4      # S33Pos
$wsProp.Range("L5").Value  = 392783.4      # S33Neg
$wsProp.Range("M5").Value  = 392783.4      # S22Pos
$wsProp.Range("N5").Value  = 392783.4      # S22Neg
$wsProp.Range("O5").Value  = 511513.2      # Z33
$wsProp.Range("P5").Value  = 511513.2      # Z22
$wsProp.Range("Q5").Value  = 99.2          # R33
$wsProp.Range("R5").Value  = 99.2          # R22

# ---------------------------------------------------------------------------
# 3) "Frame Assigns - Sect Prop": Beam B4 and Column C4 reassigned from
#    FSec1_corrosion to FSec1 (row 7 = B4, row 11 = C4)
# ---------------------------------------------------------------------------
$wsSect = $wb.Worksheets.Item("Frame Assigns - Sect Prop")
$wsSect.Range("F7").Value = "FSec1"
$wsSect.Range("F11").Value = "FSec1"

# ---------------------------------------------------------------------------
# 4) "Frame Assigns - Summary": same reassignment reflected in Analysis
#    Section / Design Section columns (F/G) for B4 (row 7) and C4 (row 11)
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Frame Assigns - Summary")
$wsSummary.Range("F7").Value = "FSec1"
$wsSummary.Range("G7").Value = "FSec1"
$wsSummary.Range("F11").Value = "FSec1"
$wsSummary.Range("G11").Value = "FSec1"

# ---------------------------------------------------------------------------
# 5) "Frame Assigns - End Len Offsets": Offset I for B2 (row 5) and Offset J
#    for B3 (row 6) corrected from 146.8 to 150
# ---------------------------------------------------------------------------
$wsOffsets = $wb.Worksheets.Item("Frame Assigns - End Len Offsets")
$wsOffsets.Range("E5").Value = 150
$wsOffsets.Range("F6").Value = 150

# ---------------------------------------------------------------------------
# 6) Minor selection/view state changes (cosmetic, matches saved file state)
# ---------------------------------------------------------------------------
$wsColConn = $wb.Worksheets.Item("Column Object Connectivity")
$wsColConn.Range("A4:E7").Select()

$wsPtConn = $wb.Worksheets.Item("Point Object Connectivity")
$wsPtConn.Range("A4:D11").Select()
